$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 336, shifting the existing rows 336:357 down to 337:358.
$ws.Rows(336).Insert()

# Populate the newly inserted row 336 with the new price-observation record.
$ws.Cells.Item(336, 1).Value = 5
$ws.Cells.Item(336, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(336, 3).Value = "Maule"
$ws.Cells.Item(336, 4).Value = 44714
$ws.Cells.Item(336, 5).Value = 7
$ws.Cells.Item(336, 6).Value = 100112032
$ws.Cells.Item(336, 7).Value = "Zapallo italiano"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 300
$ws.Cells.Item(336, 11).Value = 11000
$ws.Cells.Item(336, 12).Value = 11000
$ws.Cells.Item(336, 13).Value = 11000
$ws.Cells.Item(336, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(336, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(336, 16).Value = 220
$ws.Cells.Item(336, 17).Value = 50
$ws.Cells.Item(336, 18).Value = "Hortaliza"
